$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Cell A2 held the shared string "MODELO"; replace it with the numeric
# value captured from the drawing's updated parameters, keeping the cell
# centered both horizontally and vertically (style index 1 in the target).
$ws.Range("A2").Value = 2000049503
$ws.Range("A2").VerticalAlignment = -4108   # xlCenter
$ws.Range("A2").HorizontalAlignment = -4108 # xlCenter

# Move/record the active selection to A2 (was C6).
$ws.Range("A2").Select()
